$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 7593
$ws.Range("C2").Value = "Standard"

# Row 3
$ws.Range("B3").Value = 5639
$ws.Range("C3").Value = "Libre-Service"

# Row 4
$ws.Range("B4").Value = 3448
$ws.Range("C4").Value = "Standard"

# Row 5
$ws.Range("B5").Value = 7111

# Row 6
$ws.Range("B6").Value = 6849
$ws.Range("C6").Value = "Express"

# Row 7
$ws.Range("B7").Value = 1421

# Row 8
$ws.Range("B8").Value = 2385

# Row 9
$ws.Range("B9").Value = 1585
$ws.Range("C9").Value = "Express"

# Row 10
$ws.Range("B10").Value = 7790

# Row 11
$ws.Range("B11").Value = 4994
$ws.Range("C11").Value = "Libre-Service"
